$wb = $excel.ActiveWorkbook

# --- CreateDepartment sheet: add OrgUnit column (D) ---
$wsDept = $wb.Worksheets.Item("CreateDepartment")
$wsDept.Range("D1").Value = "OrgUnit"
$wsDept.Range("D2").Value = "India>South>Bangalore>Dell"

# --- Create sheet: add CallFlowURL / Intent / TemplateType columns (H, I, J) ---
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("H1").Value = "CallFlowURL"
$wsCreate.Range("I1").Value = "Intent"
$wsCreate.Range("I2").Value = "'89"
$wsCreate.Range("H2").Value = "Email"
$wsCreate.Range("J1").Value = "TemplateType"
$wsCreate.Range("J2").Value = "Readonly"

# --- Make "Create" the active sheet/tab with J3 selected ---
$wsCreate.Activate()
$wsCreate.Range("J3").Select()
